# impedir que um apartamento tenha mais de 2 vagas
# Marca vagas adicionais como "Pré-Selecionada" = SIM e define o
# apartamento elegível correspondente na coluna "Apartamentos Elegíveis".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vagas")

# Linha -> Apartamento elegível
$updates = @{
    5   = 201
    6   = 202
    7   = 2002
    10  = 303
    47  = 203
    48  = 204
    54  = 205
    55  = 1101
    93  = 206
    95  = 2003
    96  = 2003
    101 = 306
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = "SIM"
    $ws.Range("E$row").Value = $updates[$row]
}

# Atualiza a célula ativa selecionada para refletir a última edição.
$ws.Range("E55").Select() | Out-Null
